$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": insert a new date column ("01-dec") right before the
# existing "01-oct." column (column ED, index 134), shifting every column
# from ED onward one position to the right (ED..FH -> EE..FI).
# ---------------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")
$wsSpot.Columns.Item(134).Insert()

# Header row: label the newly inserted column.
$wsSpot.Cells.Item(1,134).Value = "01-dec"

# Data rows (2..25): fill the newly inserted column with a placeholder "-"
# like every other "missing data" cell in the sheet.
for ($r = 2; $r -le 25; $r++) {
  $wsSpot.Cells.Item($r,134).Value = "-"
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append two new daily rows after the existing last row (163).
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Cells.Item(164,1).NumberFormat = "@"
$wsGaz.Cells.Item(164,1).Value = "2025-11-29"
$wsGaz.Cells.Item(164,1).ClearFormats()
$wsGaz.Cells.Item(164,2).Value = 27.525

$wsGaz.Cells.Item(165,1).NumberFormat = "@"
$wsGaz.Cells.Item(165,1).Value = "2025-11-30"
$wsGaz.Cells.Item(165,1).ClearFormats()
$wsGaz.Cells.Item(165,2).Value = 27.525

# ---------------------------------------------------------------------------
# Sheet "CO2": append the same two new daily rows (price not yet published,
# so column B stays empty, exactly as for the last existing row 163).
# ---------------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")

$wsCO2.Cells.Item(164,1).NumberFormat = "@"
$wsCO2.Cells.Item(164,1).Value = "2025-11-29"
$wsCO2.Cells.Item(164,1).ClearFormats()

$wsCO2.Cells.Item(165,1).NumberFormat = "@"
$wsCO2.Cells.Item(165,1).Value = "2025-11-30"
$wsCO2.Cells.Item(165,1).ClearFormats()

Write-Host "edit applied"
